$wb = $excel.ActiveWorkbook

$wsTrain = $wb.Worksheets.Item("train")

$holeIds = @(
    "BRG_05_12",
    "BRG_16_06",
    "ECO_09_04",
    "BRG_16_02",
    "BRG_13_02",
    "BRG_05_04",
    "BRG_05_11",
    "BRG_01_02",
    "BRG_16_04B",
    "ECO_09_03",
    "BRG_16_01",
    "BRG_16_08",
    "BRG_05_02",
    "BRG_01_05",
    "BRG_01_07",
    "BRG_01_04",
    "BRG_16_09",
    "BRG_01_06",
    "BRG_16_04A",
    "BRG_05_15",
    "BRG_16_03",
    "ECO_09_05",
    "BRG_08_01",
    "ECO_09_01",
    "BRG_16_07",
    "BRG_01_08",
    "BRG_05_03",
    "BRG_05_13",
    "BRG_13_01",
    "BRG_01_09",
    "BRG_05_10"
)

# Match the header formatting already used by B1:M1 (bold, bordered, centered)
$wsTrain.Range("B1").Copy()
$wsTrain.Range("A1").PasteSpecial(-4122)
$wsTrain.Range("A1").Value = "hole_id"

for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = $i + 2
    $wsTrain.Cells.Item($row, 1).Value = $holeIds[$i]
}
